$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the store names between row 4 and row 5 (Manauara <-> Ponta Negra)
$ws.Range("A4").Value = "Bibi Cell Ponta Negra"
$ws.Range("A5").Value = "Bibi Cell Manauara"

# Row 2 updates
$ws.Range("O2").Value = 22226.83
$ws.Range("AG2").Value = 161641.12

# Row 3 updates
$ws.Range("O3").Value = 6484.9
$ws.Range("AG3").Value = 88631.7

# Row 4 updates (now holds what used to be row 5's daily figures)
$ws.Range("B4").Value = 2321
$ws.Range("C4").Value = 4256.58
$ws.Range("D4").Value = 3600
$ws.Range("E4").Value = 2120
$ws.Range("F4").Value = 2699.7
$ws.Range("G4").Value = 2150.81
$ws.Range("H4").Value = 2966
$ws.Range("I4").Value = 2107.11
$ws.Range("J4").Value = 2296
$ws.Range("K4").Value = 2017.01
$ws.Range("L4").Value = 6974.9
$ws.Range("M4").Value = 2702.01
$ws.Range("N4").Value = 3022.01
$ws.Range("O4").Value = 6136.45
$ws.Range("AG4").Value = 45369.58

# Row 5 updates (now holds what used to be row 4's daily figures)
$ws.Range("B5").Value = 5020
$ws.Range("C5").Value = 2670
$ws.Range("D5").Value = 4432
$ws.Range("E5").Value = 2250
$ws.Range("F5").Value = 3344.5
$ws.Range("G5").Value = 2487
$ws.Range("H5").Value = 4148
$ws.Range("I5").Value = 3419.9
$ws.Range("J5").Value = 2395.49
$ws.Range("K5").Value = 4368
$ws.Range("L5").Value = 2802
$ws.Range("M5").Value = 2668
$ws.Range("N5").Value = 3201.9
$ws.Range("O5").Value = 1561
$ws.Range("AG5").Value = 44767.79

# Row 6 updates
$ws.Range("O6").Value = 36409.18
$ws.Range("AG6").Value = 340410.19
